# Apply updated "want-to-go" counts and sold-out status to the
# "展览" and "全部类型" worksheets (both hold the same data).

$wb = $excel.ActiveWorkbook

# Map of cell -> new numeric value for column F (想去人数)
$fUpdates = @{
    "F2"  = 1149
    "F3"  = 864
    "F8"  = 2394
    "F9"  = 7789
    "F10" = 929
    "F11" = 451
    "F12" = 389
    "F13" = 160
    "F14" = 432
    "F17" = 8010
    "F19" = 1389
    "F24" = 331
    "F28" = 113
    "F30" = 428
    "F31" = 1161
    "F34" = 67
    "F35" = 85
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($cellRef in $fUpdates.Keys) {
        $ws.Range($cellRef).Value = $fUpdates[$cellRef]
    }

    # G22: status changes from "已售罄" (sold out) to "不可售" (not for sale)
    $ws.Range("G22").Value = "不可售"
}
